# Add a new "All Matches" sheet (for printing) right after "Judging Schedule".
# It lists every match, by team, in raw (Team, Match Number, Position) form.
# Manually merged from Aryan branch.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("Judging Schedule")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = "All Matches"

$matches = @(
    @(12001,4,1),
    @(12001,6,2),
    @(12001,8,3),
    @(12001,10,2),
    @(12001,12,4),
    @(12002,6,1),
    @(12002,8,2),
    @(12002,10,4),
    @(12002,12,3),
    @(12002,14,3),
    @(12003,1,1),
    @(12003,8,4),
    @(12003,10,3),
    @(12003,12,2),
    @(12003,14,2),
    @(12004,1,2),
    @(12004,3,3),
    @(12004,10,1),
    @(12004,12,1),
    @(12004,14,4),
    @(12005,1,3),
    @(12005,3,4),
    @(12005,5,2),
    @(12005,11,1),
    @(12005,13,3),
    @(12006,1,4),
    @(12006,3,2),
    @(12006,5,1),
    @(12006,7,3),
    @(12006,13,4),
    @(12007,2,3),
    @(12007,4,2),
    @(12007,6,4),
    @(12007,8,1),
    @(12007,15,4),
    @(12008,2,4),
    @(12008,4,3),
    @(12008,7,4),
    @(12008,9,2),
    @(12008,11,2),
    @(12009,2,2),
    @(12009,4,4),
    @(12009,6,3),
    @(12009,9,1),
    @(12009,11,4)
)

for ($i = 0; $i -lt $matches.Count; $i++) {
    $row = $matches[$i]
    $ws.Cells.Item($i + 1, 1).Value = $row[0]
    $ws.Cells.Item($i + 1, 2).Value = $row[1]
    $ws.Cells.Item($i + 1, 3).Value = $row[2]
}

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("Match Schedule").Activate()
